$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# C3 / D3 used to hold the shared strings "name" / "surename" (likely populated
# from the Person array that TicketWriter pulls the booking's name/surname
# from). Clear them out to empty cells.
$ws.Range("C3").Value = ""
$ws.Range("C3").Borders.LineStyle = -4142

$ws.Range("D3").Value = ""
$ws.Range("D3").Borders.LineStyle = -4142

# H3 / C8 / D8 used to be populated with hard-coded placeholder numbers
# (69 / 7331 / 1337); now wired up to the (currently empty) array lists, so
# they resolve to 0.
$ws.Range("H3").Value = 0
$ws.Range("C8").Value = 0
$ws.Range("D8").Value = 0
